# CCMA.xlsx — "Update to match all regions"
#
# The sheet builds a "AAAAMM" (year+month) tag from TODAY() via
# TEXT(date,"AAAAMM"). That custom format code is Spanish-locale-specific
# and doesn't resolve the same way everywhere, so it's replaced with the
# locale-independent YEAR(date)&TEXT(MONTH(date),"00").
#
# The formula lives in column J as three shared-formula groups (anchors
# J2, J3 [[spanning J3:J66]] and J67 [[spanning J67:J69]]); each group is
# rewritten as a whole range so Excel keeps it as one shared formula
# instead of silently exploding it into per-cell formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 — single (non-shared) formula cell.
$ws.Range("J2").Formula = '=CONCATENATE(TEXT(A2,"0")," - ","CCMA - ",YEAR(F2)&TEXT(MONTH(F2),"00")," - ",SUBSTITUTE(D2,"-","")," - ",B2)'

# Rows 3-66 — shared formula group anchored at J3.
$ws.Range("J3:J66").Formula = '=CONCATENATE(TEXT(A3,"0")," - ","CCMA - ",YEAR(F3)&TEXT(MONTH(F3),"00")," - ",SUBSTITUTE(D3,"-","")," - ",B3)'

# Rows 67-69 — shared formula group anchored at J67.
$ws.Range("J67:J69").Formula = '=CONCATENATE(TEXT(A67,"0")," - ","CCMA - ",YEAR(F67)&TEXT(MONTH(F67),"00")," - ",SUBSTITUTE(D67,"-","")," - ",B67)'

# The saved view had the cursor left on D5; re-select the frozen pane's
# home cell (A2) to match the refreshed snapshot.
$ws.Range("A2").Select() | Out-Null
